# Add a new closing "Köszönjük a figyelmet" (Thank you for your attention)
# slide at the end of the deck, matching the look of the other
# "Cím és tartalom" (Title and Content) slides already in the file.

$p = $ppt.ActivePresentation

# Duplicate the last existing slide so the new slide inherits the same
# slide layout, placeholder naming/ids, color map override, etc. as the
# rest of the deck.
$lastSlide = $p.Slides.Item($p.Slides.Count)
$lastSlide.Duplicate() | Out-Null

# The duplicate was inserted right after $lastSlide, i.e. it is now the
# last slide in the deck.
$newSlide = $p.Slides.Item($p.Slides.Count)

# Title placeholder -> new closing title.
$titleShape = $newSlide.Shapes.Item(1)
$titleShape.TextFrame.TextRange.Text = "Köszönjük a figyelmet"

# Content placeholder -> clear out the copied bullet list, leaving it
# empty like on the freshly-created slide.
$contentShape = $newSlide.Shapes.Item(2)
$contentShape.TextFrame.TextRange.Text = ""
